$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to an exact text value, preserving the original
# (unstyled) cell format even when the text looks like a number.
function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "242.57"
Set-TextValue $ws.Range("D4") "5.287"
Set-TextValue $ws.Range("D5") "0.05626"
Set-TextValue $ws.Range("D6") "3.379"
Set-TextValue $ws.Range("D7") "6.372"
Set-TextValue $ws.Range("D8") "0.8070"
Set-TextValue $ws.Range("D9") "0.9574"
Set-TextValue $ws.Range("D10") "0.1427"
Set-TextValue $ws.Range("D11") "0.07446"
Set-TextValue $ws.Range("D12") "0.03228"
Set-TextValue $ws.Range("D13") "0.03067"
Set-TextValue $ws.Range("D14") "0.09275"
Set-TextValue $ws.Range("D15") "3.568"
Set-TextValue $ws.Range("D16") "0.001651"
Set-TextValue $ws.Range("D17") "0.04709"
Set-TextValue $ws.Range("D18") "0.0005828"
Set-TextValue $ws.Range("D19") "0.006356"
Set-TextValue $ws.Range("D20") "0.004976"
Set-TextValue $ws.Range("D21") "0.001042"
Set-TextValue $ws.Range("D22") "0.0001502"
Set-TextValue $ws.Range("D23") "0.0003104"
Set-TextValue $ws.Range("D24") "3.770"
Set-TextValue $ws.Range("D25") "2.094"
Set-TextValue $ws.Range("D27") "0.1275"
Set-TextValue $ws.Range("D40") "0.03916"
Set-TextValue $ws.Range("D41") "0.006949"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D42") "0.1035"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D43") "0.002913"
$ws.Range("E43").Value = "42CEJICEJI"
Set-TextValue $ws.Range("D44") "0.007480"
Set-TextValue $ws.Range("D45") "0.00005940"
Set-TextValue $ws.Range("D47") "0.0005507"
Set-TextValue $ws.Range("D48") "0.6834"
Set-TextValue $ws.Range("D49") "0.05851"
Set-TextValue $ws.Range("D50") "0.00002103"
